$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rule R20's "Integer min" (C10) was restored from 18 back to 1.
$ws.Range("C10").Value = 1
